$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 32×82=2624 -> 11×90=990
$cell = $t.Cell(1, 1)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "11×90=990"

# 89×82=7298 -> 67×80=5360
$cell = $t.Cell(1, 2)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "67×80=5360"

# 24×65=1560 -> 78×85=6630
$cell = $t.Cell(1, 3)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "78×85=6630"

# 40×77=3080 -> 51×43=2193
$cell = $t.Cell(1, 4)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "51×43=2193"

# 87×38=3306 -> 21×58=1218
$cell = $t.Cell(1, 5)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "21×58=1218"

# 23×55=1265 -> 93×46=4278
$cell = $t.Cell(5, 1)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "93×46=4278"

# 40×24=960 -> 64×57=3648
$cell = $t.Cell(5, 2)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "64×57=3648"

# 84×36=3024 -> 92×97=8924
$cell = $t.Cell(5, 3)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "92×97=8924"

# 49×62=3038 -> 14×15=210
$cell = $t.Cell(5, 4)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "14×15=210"

# 34×61=2074 -> 64×68=4352
$cell = $t.Cell(5, 5)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "64×68=4352"

# 61×84=5124 -> 63×52=3276
$cell = $t.Cell(10, 1)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "63×52=3276"

# 46×43=1978 -> 33×43=1419
$cell = $t.Cell(10, 2)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "33×43=1419"

# 18×74=1332 -> 40×25=1000
$cell = $t.Cell(10, 3)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "40×25=1000"

# 57×42=2394 -> 85×19=1615
$cell = $t.Cell(10, 4)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "85×19=1615"

# 39×87=3393 -> 63×11=693
$cell = $t.Cell(10, 5)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "63×11=693"

# 89×88=7832 -> 30×90=2700
$cell = $t.Cell(15, 1)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "30×90=2700"

# 72×91=6552 -> 49×33=1617
$cell = $t.Cell(15, 2)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "49×33=1617"

# 23×71=1633 -> 34×61=2074
$cell = $t.Cell(15, 3)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "34×61=2074"

# 33×78=2574 -> 92×65=5980
$cell = $t.Cell(15, 4)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "92×65=5980"

# 49×16=784 -> 78×61=4758
$cell = $t.Cell(15, 5)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "78×61=4758"

# 67×43=2881 -> 37×41=1517
$cell = $t.Cell(20, 1)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "37×41=1517"

# 37×57=2109 -> 56×86=4816
$cell = $t.Cell(20, 2)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "56×86=4816"

# 60×35=2100 -> 38×13=494
$cell = $t.Cell(20, 3)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "38×13=494"

# 77×24=1848 -> 45×66=2970
$cell = $t.Cell(20, 4)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "45×66=2970"

# 28×35=980 -> 93×23=2139
$cell = $t.Cell(20, 5)
$r = $cell.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "93×23=2139"
